# Auto-generated edit script for A 54142-2021.docx

function New-Para {
    param($doc)
    $r = $doc.Content
    $r.Collapse(0)
    $insertPoint = $r.End
    $r.InsertParagraphAfter()
    return $insertPoint
}

function Add-PlainRun {
    param($doc, $pos, $text)
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    return $r.End
}

function Add-ItalicRun {
    param($doc, $pos, $text)
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    $endPos = $r.End
    $fr = $doc.Range($pos, $endPos)
    $fr.Font.Italic = 1
    return $endPos
}

$d = $word.ActiveDocument

# Paragraph 1 (style=Heading1)
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Knärot – ekologi samt krav på livsmiljön'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 1"

# Paragraph 2
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 3
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Samuel Johnsons doktorsavhandling '
$pos = Add-ItalicRun $d $pos '“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“'
$pos = Add-PlainRun $d $pos ' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: '
$pos = Add-ItalicRun $d $pos '“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” '
$pos = Add-PlainRun $d $pos 'Vidare '
$pos = Add-ItalicRun $d $pos '“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 4
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: '
$pos = Add-ItalicRun $d $pos '“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 5
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 6
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 7 (style=Heading2)
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Referenser - knärot'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Heading 2"

# Paragraph 8
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'de Graaf M & Roberts M.R., 2009. '
$pos = Add-ItalicRun $d $pos 'Short-term response of the herbaceous layer within leave patches after harvest. '
$pos = Add-PlainRun $d $pos 'Forest Ecology and Management 257, 1014-1025'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 9
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. '
$pos = Add-ItalicRun $d $pos 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. '
$pos = Add-PlainRun $d $pos 'Ecological Applications, 22, 2049-2064 '
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 10
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. '
$pos = Add-ItalicRun $d $pos 'Interactive effects of drought and edge exposure on old-growth forest understory species. '
$pos = Add-PlainRun $d $pos 'Landscape Ecology, 37, sid 1839-1853'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 11
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. '
$pos = Add-ItalicRun $d $pos 'Biological legacies buffer local species extinction after logging. '
$pos = Add-PlainRun $d $pos 'Journal of Applied Ecology. 51, 53-62.'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 12
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'Skogsstyrelsen, 2022. '
$pos = Add-ItalicRun $d $pos 'Vägledning för hänsyn till knärot. '
$pos = Add-PlainRun $d $pos 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Paragraph 13
$pos = New-Para $d
$pos = Add-PlainRun $d $pos 'SLU Artdatabanken, 2021. '
$pos = Add-ItalicRun $d $pos 'Artfaktablad. Naturvård – artfakta. '
$pos = Add-PlainRun $d $pos 'SLU Artdatabanken, Uppsala '
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"

# Update date in first-page header
$sec = $d.Sections.Item(1)
$headers = $sec.Headers
$firstPageHeader = $headers.Item(2)
$found = $firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
Write-Host "Date replaced:" $found

Write-Host "Total paragraphs:" $d.Paragraphs.Count
